$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the former FAPs-target row (with refreshed TPM-derived values)
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08844
$ws.Range("H2").Value = 0.26532
$ws.Range("I2").Value = 0.939488472392877
$ws.Range("J2").Value = 0.9394884723928769
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.26202
$ws.Range("N2").Value = 0.78606
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0231730488
$ws.Range("R2").Value = 0.2085574392
$ws.Range("S2").Value = 0.939488472392877
$ws.Range("T2").Value = 0.9394884723928769

# Row 3 becomes the former MuSCs-target row (with refreshed TPM-derived values)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.005696333333333334
$ws.Range("H3").Value = 0.017089
$ws.Range("I3").Value = 0.06051152760712301
$ws.Range("J3").Value = 0.060511527607123
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.26202
$ws.Range("N3").Value = 0.78606
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.00149255326
$ws.Range("R3").Value = 0.01343297934
$ws.Range("S3").Value = 0.06051152760712301
$ws.Range("T3").Value = 0.060511527607123

# Remove the now-obsolete ECs / MuSCs(row4) / Resolving-Mac (row5) source rows
$ws.Range("A4:T5").Delete()
